$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "22.339.06"
$ws.Range("E2").Value = "  -0.10%  "

Set-TextValue "D3" "1.563.03"
$ws.Range("E3").Value = "  +0.00%  "

Set-TextValue "D4" "1.007"
$ws.Range("E4").Value = "  +0.52%  "

Set-TextValue "D5" "1.006"
$ws.Range("E5").Value = "  +0.47%  "

Set-TextValue "D6" "289.14"
$ws.Range("E6").Value = "  -0.45%  "

Set-TextValue "D7" "0.3737"
$ws.Range("E7").Value = "  +0.63%  "

Set-TextValue "D8" "49.45"
$ws.Range("E8").Value = "  +0.85%  "

Set-TextValue "D9" "0.3361"
$ws.Range("E9").Value = "  -0.94%  "

Set-TextValue "D10" "0.07461"
$ws.Range("E10").Value = "  -2.34%  "

Set-TextValue "D11" "1.115"
$ws.Range("E11").Value = "  -4.31%  "

Set-TextValue "D12" "1.007"
$ws.Range("E12").Value = "  +0.45%  "

Set-TextValue "D13" "20.69"
$ws.Range("E13").Value = "  -3.56%  "

Set-TextValue "D14" "5.856"
$ws.Range("E14").Value = "  -3.20%  "

Set-TextValue "D15" "6.849"
$ws.Range("E15").Value = "  -0.98%  "

Set-TextValue "D16" "1.564.52"
$ws.Range("E16").Value = "  -0.01%  "

Set-TextValue "D17" "0.00001102"
$ws.Range("E17").Value = "  -2.18%  "

Set-TextValue "D18" "89.02"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("E19").Value = "  -0.15%  "

Set-TextValue "D20" "1.006"
$ws.Range("E20").Value = "  +0.43%  "

Set-TextValue "D21" "6.155"
$ws.Range("E21").Value = "  -1.21%  "

Set-TextValue "D22" "16.21"
$ws.Range("E22").Value = "  -1.99%  "

Set-TextValue "D23" "11.81"
$ws.Range("E23").Value = "  -1.73%  "

Set-TextValue "D24" "22.327.38"
$ws.Range("E24").Value = "  -0.16%  "

Set-TextValue "D25" "2.373"
$ws.Range("E25").Value = "  -1.22%  "

Set-TextValue "D26" "2.597"
$ws.Range("E26").Value = "  -7.63%  "

Set-TextValue "D27" "19.89"
$ws.Range("E27").Value = "  -1.39%  "

Set-TextValue "D28" "147.44"
$ws.Range("E28").Value = "  +1.54%  "

Set-TextValue "D29" "5.007"
$ws.Range("E29").Value = "  +0.46%  "

Set-TextValue "D30" "124.39"
$ws.Range("E30").Value = "  -0.76%  "

Set-TextValue "D31" "1.738.67"
$ws.Range("E31").Value = "  -0.07%  "

Set-TextValue "D32" "2.020"
$ws.Range("E32").Value = "  +0.55%  "

Set-TextValue "D33" "0.9770"
$ws.Range("E33").Value = "  -2.63%  "

Set-TextValue "D34" "5.869"
$ws.Range("E34").Value = "  -5.33%  "

Set-TextValue "D35" "9.739"
$ws.Range("E35").Value = "  -2.87%  "

$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D36" "1.406"
$ws.Range("E36").Value = "  +9.11%  "

$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D37" "0.08428"
$ws.Range("E37").Value = "  -0.56%  "

Set-TextValue "D38" "0.02440"
$ws.Range("E38").Value = "  -3.46%  "

Set-TextValue "D39" "0.2248"
$ws.Range("E39").Value = "  -3.23%  "

Set-TextValue "D40" "0.06367"
$ws.Range("E40").Value = "  -0.35%  "

Set-TextValue "D41" "5.321"
$ws.Range("E41").Value = "  -3.55%  "

Set-TextValue "D42" "0.6210"
$ws.Range("E42").Value = "  -2.09%  "

Set-TextValue "D43" "10.87"
$ws.Range("E43").Value = "  -6.99%  "

Set-TextValue "D44" "1.005"
$ws.Range("E44").Value = "  +0.40%  "

Set-TextValue "D45" "13.90"
$ws.Range("E45").Value = "  -1.47%  "

Set-TextValue "D46" "3.787"
$ws.Range("E46").Value = "  +0.70%  "

Set-TextValue "D47" "0.5737"
$ws.Range("E47").Value = "  -3.91%  "

Set-TextValue "D48" "2.034"
$ws.Range("E48").Value = "  -2.83%  "

Set-TextValue "D49" "1.236"
$ws.Range("E49").Value = "  -2.25%  "

Set-TextValue "D50" "123.79"
$ws.Range("E50").Value = "  -0.58%  "

Set-TextValue "D51" "0.07287"
$ws.Range("E51").Value = "  +0.28%  "
